$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 81: height changed from 33.6 to 16.8 (no content change) ---
$ws.Rows.Item(81).RowHeight = 16.8

# --- Copy style templates once: style1/style3 pair (from row 187), style1/style6 pair (from row 170) ---
$ws.Range("A187:B187").Copy()
$styleRows = @(188,189,190,191,192,193,194,195,196,197,199,200,201,202,203,204,205,206,207,208,209,210,211,212,213,214,215)
foreach ($r in $styleRows) {
    $ws.Range("A" + $r + ":B" + $r).PasteSpecial(-4122)
}

$ws.Range("A170:B170").Copy()
$ws.Range("A198:B198").PasteSpecial(-4122)

# --- Fill in values. Row 193 and 205 get column B written before column A ---
# --- to reproduce the original shared-string insertion order.            ---
$ws.Range("A188").Value = "Breigleb, Johann Christian"
$ws.Range("B188").Value = "http://viaf.org/viaf/35210135"
$ws.Range("A189").Value = "Brix, Ernst Julius"
$ws.Range("B189").Value = "http://viaf.org/viaf/35210587"
$ws.Range("A190").Value = "Brockhard, Michael"
$ws.Range("B190").Value = "http://viaf.org/viaf/191166404"
$ws.Range("A191").Value = "Brodribb, William Jackson"
$ws.Range("B191").Value = "http://viaf.org/viaf/76716202"
$ws.Range("A192").Value = "Church, Alfred John"
$ws.Range("B192").Value = "http://viaf.org/viaf/45233331"
$ws.Range("B193").Value = "http://viaf.org/viaf/164447798"
$ws.Range("A193").Value = "Brohm, Karl Friedr. Aug."
$ws.Range("A194").Value = "Brossaeus, C."
$ws.Range("B194").Value = "http://viaf.org/viaf/59079671"
$ws.Range("A195").Value = "Browning, Oscar"
$ws.Range("B195").Value = "http://viaf.org/viaf/130145857831723020397"
$ws.Range("A196").Value = "Bruder, Carl Hermann"
$ws.Range("B196").Value = "http://viaf.org/viaf/18273247"
$ws.Range("A197").Value = "Brunck, Richard Francois-Phillipe"
$ws.Range("B197").Value = "http://viaf.org/viaf/82656304"
$ws.Range("A198").Value = "Bruno, Agostino"
$ws.Range("B198").Value = "http://viaf.org/viaf/22515703"
$ws.Range("A199").Value = "Bruns, Paul Jakob"
$ws.Range("B199").Value = "http://viaf.org/viaf/17983318"
$ws.Range("A200").Value = "Bryce, Archilbald Hamilton"
$ws.Range("B200").Value = "http://viaf.org/viaf/89610712"
$ws.Range("A201").Value = "Buckley, Theodore Alors"
$ws.Range("B201").Value = "http://viaf.org/viaf/45879798"
$ws.Range("A202").Value = "Budai, Ezsaias"
$ws.Range("B202").Value = "http://viaf.org/viaf/39386168"
$ws.Range("A203").Value = "Buenemann, Johann Ludolf"
$ws.Range("B203").Value = "http://viaf.org/viaf/780439"
$ws.Range("A204").Value = "Burman, Kasper"
$ws.Range("B204").Value = "http://viaf.org/viaf/88788037"
$ws.Range("B205").Value = "http://viaf.org/viaf/47055143"
$ws.Range("A205").Value = "Reiske, Johann Jacob"
$ws.Range("A206").Value = "Franz, Johann Georg Friedrich"
$ws.Range("B206").Value = "http://viaf.org/viaf/12619940"
$ws.Range("A207").Value = "Beck, Christian Daniel"
$ws.Range("B207").Value = "http://viaf.org/viaf/57357633"
$ws.Range("A208").Value = "Matthiae, F.C."
$ws.Range("B208").Value = "http://viaf.org/viaf/69692810"
$ws.Range("A209").Value = "Heinsius, Daniel"
$ws.Range("B209").Value = "http://viaf.org/viaf/56635500"
$ws.Range("A210").Value = "Heinsius, Nicolaas"
$ws.Range("B210").Value = "http://viaf.org/viaf/95302161"
$ws.Range("A211").Value = "Ciofano, Ercole"
$ws.Range("B211").Value = "http://viaf.org/viaf/45056369"
$ws.Range("A212").Value = "Moltzer, Jakob"
$ws.Range("B212").Value = "http://viaf.org/viaf/56633987"
$ws.Range("A213").Value = "Burnouf, Emile"
$ws.Range("B213").Value = "http://viaf.org/viaf/56629753"
$ws.Range("A214").Value = "Burton, Edward"
$ws.Range("B214").Value = "http://viaf.org/viaf/75042825"
$ws.Range("A215").Value = "Buttmann, Phillip"
$ws.Range("B215").Value = "http://viaf.org/viaf/39611820"

# --- Row heights for the new rows (all 16.8 except row 198, which keeps default) ---
$heightRows = @(188,189,190,191,192,193,194,195,196,197,199,200,201,202,203,204,205,206,207,208,209,210,211,212,213,214,215)
foreach ($r in $heightRows) {
    $ws.Rows.Item($r).RowHeight = 16.8
}

# --- Final selection: Excel ends up with A216 selected after data entry ---
$ws.Range("A216").Select()
